$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("G2").Value = 19.434592
$ws.Range("H2").Value = 58.303776
$ws.Range("I2").Value = 0.1244167820899015
$ws.Range("J2").Value = 0.1244167820899015
$ws.Range("K2").Value = 3.0
$ws.Range("M2").Value = 0.7597586666666668
$ws.Range("N2").Value = 2.279276
$ws.Range("O2").Value = 0.3177111952105157
$ws.Range("P2").Value = 0.3177111952105158
$ws.Range("Q2").Value = 14.76559970513067
$ws.Range("R2").Value = 132.890397346176
$ws.Range("S2").Value = 0.0395286045420289
$ws.Range("T2").Value = 0.0395286045420289
$ws.Range("E3").Value = 3.0
$ws.Range("G3").Value = 19.434592
$ws.Range("H3").Value = 58.303776
$ws.Range("I3").Value = 0.1244167820899015
$ws.Range("J3").Value = 0.1244167820899015
$ws.Range("K3").Value = 3.0
$ws.Range("M3").Value = 0.6247346666666667
$ws.Range("N3").Value = 1.874204
$ws.Range("O3").Value = 0.2612476913319534
$ws.Range("P3").Value = 0.2612476913319534
$ws.Range("Q3").Value = 12.14146335492267
$ws.Range("R3").Value = 109.273170194304
$ws.Range("S3").Value = 0.0325035970839375
$ws.Range("T3").Value = 0.0325035970839375
$ws.Range("E4").Value = 3.0
$ws.Range("G4").Value = 19.434592
$ws.Range("H4").Value = 58.303776
$ws.Range("I4").Value = 0.1244167820899015
$ws.Range("J4").Value = 0.1244167820899015
$ws.Range("K4").Value = 3.0
$ws.Range("M4").Value = 0.3417453333333333
$ws.Range("N4").Value = 1.025236
$ws.Range("O4").Value = 0.1429089565865864
$ws.Range("P4").Value = 0.1429089565865864
$ws.Range("Q4").Value = 6.641681121237333
$ws.Range("R4").Value = 59.775130091136
$ws.Range("S4").Value = 0.01778027251032852
$ws.Range("T4").Value = 0.01778027251032852
$ws.Range("E5").Value = 3.0
$ws.Range("G5").Value = 19.434592
$ws.Range("H5").Value = 58.303776
$ws.Range("I5").Value = 0.1244167820899015
$ws.Range("J5").Value = 0.1244167820899015
$ws.Range("K5").Value = 3.0
$ws.Range("M5").Value = 0.6651113333333333
$ws.Range("N5").Value = 1.995334
$ws.Range("O5").Value = 0.2781321568709446
$ws.Range("P5").Value = 0.2781321568709446
$ws.Range("Q5").Value = 12.92616739790933
$ws.Range("R5").Value = 116.335506581184
$ws.Range("S5").Value = 0.03460430795360662
$ws.Range("T5").Value = 0.03460430795360662
$ws.Range("E6").Value = 3.0
$ws.Range("G6").Value = 48.891945
$ws.Range("H6").Value = 146.675835
$ws.Range("I6").Value = 0.3129974875220664
$ws.Range("J6").Value = 0.3129974875220664
$ws.Range("K6").Value = 3.0
$ws.Range("M6").Value = 0.7597586666666668
$ws.Range("N6").Value = 2.279276
$ws.Range("O6").Value = 0.3177111952105157
$ws.Range("P6").Value = 0.3177111952105158
$ws.Range("Q6").Value = 37.14607894394
$ws.Range("R6").Value = 334.3147104954601
$ws.Range("S6").Value = 0.0994428058585242
$ws.Range("T6").Value = 0.09944280585852419
$ws.Range("E7").Value = 3.0
$ws.Range("G7").Value = 48.891945
$ws.Range("H7").Value = 146.675835
$ws.Range("I7").Value = 0.3129974875220664
$ws.Range("J7").Value = 0.3129974875220664
$ws.Range("K7").Value = 3.0
$ws.Range("M7").Value = 0.6247346666666667
$ws.Range("N7").Value = 1.874204
$ws.Range("O7").Value = 0.2612476913319534
$ws.Range("P7").Value = 0.2612476913319534
$ws.Range("Q7").Value = 30.54449296226
$ws.Range("R7").Value = 274.90043666034
$ws.Range("S7").Value = 0.08176987100784172
$ws.Range("T7").Value = 0.08176987100784171
$ws.Range("E8").Value = 3.0
$ws.Range("G8").Value = 48.891945
$ws.Range("H8").Value = 146.675835
$ws.Range("I8").Value = 0.3129974875220664
$ws.Range("J8").Value = 0.3129974875220664
$ws.Range("K8").Value = 3.0
$ws.Range("M8").Value = 0.3417453333333333
$ws.Range("N8").Value = 1.025236
$ws.Range("O8").Value = 0.1429089565865864
$ws.Range("P8").Value = 0.1429089565865864
$ws.Range("Q8").Value = 16.70859404134
$ws.Range("R8").Value = 150.37734637206
$ws.Range("S8").Value = 0.0447301443560016
$ws.Range("T8").Value = 0.0447301443560016
$ws.Range("E9").Value = 3.0
$ws.Range("G9").Value = 48.891945
$ws.Range("H9").Value = 146.675835
$ws.Range("I9").Value = 0.3129974875220664
$ws.Range("J9").Value = 0.3129974875220664
$ws.Range("K9").Value = 3.0
$ws.Range("M9").Value = 0.6651113333333333
$ws.Range("N9").Value = 1.995334
$ws.Range("O9").Value = 0.2781321568709446
$ws.Range("P9").Value = 0.2781321568709446
$ws.Range("Q9").Value = 32.51858672821
$ws.Range("R9").Value = 292.66728055389
$ws.Range("S9").Value = 0.08705466629969888
$ws.Range("T9").Value = 0.08705466629969887
$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 65.19353
$ws.Range("H10").Value = 195.58059
$ws.Range("I10").Value = 0.4173573191390618
$ws.Range("J10").Value = 0.4173573191390618
$ws.Range("K10").Value = 3.0
$ws.Range("M10").Value = 0.7597586666666668
$ws.Range("N10").Value = 2.279276
$ws.Range("O10").Value = 0.3177111952105157
$ws.Range("P10").Value = 0.3177111952105158
$ws.Range("Q10").Value = 49.53134942809334
$ws.Range("R10").Value = 445.7821448528401
$ws.Range("S10").Value = 0.132599092693528
$ws.Range("T10").Value = 0.132599092693528
$ws.Range("E11").Value = 3.0
$ws.Range("G11").Value = 65.19353
$ws.Range("H11").Value = 195.58059
$ws.Range("I11").Value = 0.4173573191390618
$ws.Range("J11").Value = 0.4173573191390618
$ws.Range("K11").Value = 3.0
$ws.Range("M11").Value = 0.6247346666666667
$ws.Range("N11").Value = 1.874204
$ws.Range("O11").Value = 0.2612476913319534
$ws.Range("P11").Value = 0.2612476913319534
$ws.Range("Q11").Value = 40.72865823337333
$ws.Range("R11").Value = 366.55792410036
$ws.Range("S11").Value = 0.1090336360855732
$ws.Range("T11").Value = 0.1090336360855732
$ws.Range("E12").Value = 3.0
$ws.Range("G12").Value = 65.19353
$ws.Range("H12").Value = 195.58059
$ws.Range("I12").Value = 0.4173573191390618
$ws.Range("J12").Value = 0.4173573191390618
$ws.Range("K12").Value = 3.0
$ws.Range("M12").Value = 0.3417453333333333
$ws.Range("N12").Value = 1.025236
$ws.Range("O12").Value = 0.1429089565865864
$ws.Range("P12").Value = 0.1429089565865864
$ws.Range("Q12").Value = 22.27958464102667
$ws.Range("R12").Value = 200.51626176924
$ws.Range("S12").Value = 0.05964409900193825
$ws.Range("T12").Value = 0.05964409900193827
$ws.Range("E13").Value = 3.0
$ws.Range("G13").Value = 65.19353
$ws.Range("H13").Value = 195.58059
$ws.Range("I13").Value = 0.4173573191390618
$ws.Range("J13").Value = 0.4173573191390618
$ws.Range("K13").Value = 3.0
$ws.Range("M13").Value = 0.6651113333333333
$ws.Range("N13").Value = 1.995334
$ws.Range("O13").Value = 0.2781321568709446
$ws.Range("P13").Value = 0.2781321568709446
$ws.Range("Q13").Value = 43.36095566300666
$ws.Range("R13").Value = 390.2486009670599
$ws.Range("S13").Value = 0.1160804913580224
$ws.Range("T13").Value = 0.1160804913580224
$ws.Range("E14").Value = 3.0
$ws.Range("G14").Value = 22.685484
$ws.Range("H14").Value = 68.05645200000001
$ws.Range("I14").Value = 0.1452284112489703
$ws.Range("J14").Value = 0.1452284112489703
$ws.Range("K14").Value = 3.0
$ws.Range("M14").Value = 0.7597586666666668
$ws.Range("N14").Value = 2.279276
$ws.Range("O14").Value = 0.3177111952105157
$ws.Range("P14").Value = 0.3177111952105158
$ws.Range("Q14").Value = 17.23549307652801
$ws.Range("R14").Value = 155.119437688752
$ws.Range("S14").Value = 0.04614069211643466
$ws.Range("T14").Value = 0.04614069211643466
$ws.Range("E15").Value = 3.0
$ws.Range("G15").Value = 22.685484
$ws.Range("H15").Value = 68.05645200000001
$ws.Range("I15").Value = 0.1452284112489703
$ws.Range("J15").Value = 0.1452284112489703
$ws.Range("K15").Value = 3.0
$ws.Range("M15").Value = 0.6247346666666667
$ws.Range("N15").Value = 1.874204
$ws.Range("O15").Value = 0.2612476913319534
$ws.Range("P15").Value = 0.2612476913319534
$ws.Range("Q15").Value = 14.172408284912
$ws.Range("R15").Value = 127.551674564208
$ws.Range("S15").Value = 0.03794058715460098
$ws.Range("T15").Value = 0.03794058715460098
$ws.Range("E16").Value = 3.0
$ws.Range("G16").Value = 22.685484
$ws.Range("H16").Value = 68.05645200000001
$ws.Range("I16").Value = 0.1452284112489703
$ws.Range("J16").Value = 0.1452284112489703
$ws.Range("K16").Value = 3.0
$ws.Range("M16").Value = 0.3417453333333333
$ws.Range("N16").Value = 1.025236
$ws.Range("O16").Value = 0.1429089565865864
$ws.Range("P16").Value = 0.1429089565865864
$ws.Range("Q16").Value = 7.752658291408001
$ws.Range("R16").Value = 69.77392462267201
$ws.Range("S16").Value = 0.020754440718318
$ws.Range("T16").Value = 0.02075444071831801
$ws.Range("E17").Value = 3.0
$ws.Range("G17").Value = 22.685484
$ws.Range("H17").Value = 68.05645200000001
$ws.Range("I17").Value = 0.1452284112489703
$ws.Range("J17").Value = 0.1452284112489703
$ws.Range("K17").Value = 3.0
$ws.Range("M17").Value = 0.6651113333333333
$ws.Range("N17").Value = 1.995334
$ws.Range("O17").Value = 0.2781321568709446
$ws.Range("P17").Value = 0.2781321568709446
$ws.Range("Q17").Value = 15.088372510552
$ws.Range("R17").Value = 135.795352594968
$ws.Range("S17").Value = 0.04039269125961666
$ws.Range("T17").Value = 0.04039269125961666
